# Update "想去人数" (F column) counts and one "最低票价" (G13 on 演出)
# that became unavailable ("不可售") across the four sheets of the workbook,
# matching the data refresh captured in the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# ---- 展览 (Exhibitions) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 271
$ws.Range("F3").Value  = 646
$ws.Range("F6").Value  = 2835
$ws.Range("F7").Value  = 534
$ws.Range("F8").Value  = 57
$ws.Range("F10").Value = 584
$ws.Range("F11").Value = 27
$ws.Range("F12").Value = 336
$ws.Range("F14").Value = 6011
$ws.Range("F15").Value = 634
$ws.Range("F16").Value = 1050
$ws.Range("F17").Value = 15
$ws.Range("F18").Value = 243
$ws.Range("F20").Value = 89
$ws.Range("F21").Value = 550
$ws.Range("F22").Value = 3
$ws.Range("F23").Value = 42
$ws.Range("F25").Value = 129
$ws.Range("F26").Value = 1324
$ws.Range("F29").Value = 51
$ws.Range("F30").Value = 2068
$ws.Range("F31").Value = 186
$ws.Range("F32").Value = 7
$ws.Range("F33").Value = 359
$ws.Range("F35").Value = 3318

# ---- 演出 (Performances) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value  = 24
$ws.Range("F8").Value  = 93
$ws.Range("G13").Value = "不可售"
$ws.Range("F24").Value = 4058
$ws.Range("F28").Value = 151
$ws.Range("F30").Value = 70

# ---- 本地生活 (Local life) ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value  = 87
$ws.Range("F5").Value  = 2603
$ws.Range("F8").Value  = 1498
$ws.Range("F12").Value = 661

# ---- 全部类型 (All types) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 87
$ws.Range("F4").Value  = 2603
$ws.Range("F6").Value  = 1498
$ws.Range("F9").Value  = 271
$ws.Range("F10").Value = 646
$ws.Range("F11").Value = 2835
$ws.Range("F12").Value = 57
$ws.Range("F13").Value = 661
$ws.Range("F14").Value = 584
$ws.Range("F15").Value = 93
$ws.Range("F16").Value = 27
$ws.Range("F17").Value = 336
$ws.Range("F19").Value = 6011
$ws.Range("F21").Value = 634
$ws.Range("F22").Value = 1050
$ws.Range("F23").Value = 15
$ws.Range("F24").Value = 243
$ws.Range("F26").Value = 89
$ws.Range("F27").Value = 550
$ws.Range("F37").Value = 151
$ws.Range("F40").Value = 51
$ws.Range("F41").Value = 70
$ws.Range("F43").Value = 2068
$ws.Range("F46").Value = 186
$ws.Range("F47").Value = 359
$ws.Range("F49").Value = 3318

Write-Host "Applied 61 cell updates across 4 sheets."
